$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. Remove the existing "_GoBack" bookmark (it currently sits at the end of
#    paragraph 5, right after "...browser window."). It will be re-added
#    inside paragraph 1 below, matching the target layout.
# ---------------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# ---------------------------------------------------------------------------
# 2. Paragraph 1: "1. Run <p> `mvn clean package`"
#    becomes "1. From within the UI Top Level Directory, Run `mvn clean package`"
# ---------------------------------------------------------------------------
$p1 = $d.Paragraphs(1).Range
$p1.Find.Execute("Run <p> ``", $true, $false, $false, $false, $false, $true, 1, $false, "From within the UI Top Level Directory, Run ``", 2)

# ---------------------------------------------------------------------------
# 3. Re-insert the "_GoBack" bookmark inside paragraph 1, right after the new
#    "From within the UI Top Level Directory, " text and before "Run `".
# ---------------------------------------------------------------------------
$p1 = $d.Paragraphs(1).Range
$anchor = $p1.Duplicate
$anchor.Find.Execute("From within the UI Top Level Directory, ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$bmRange = $d.Range($anchor.End, $anchor.End)
$d.Bookmarks.Add("_GoBack", $bmRange)

# ---------------------------------------------------------------------------
# 4. Last paragraph (NOTE): drop the leading "<p> <b>" markup remnants.
# ---------------------------------------------------------------------------
$lastPar = $d.Paragraphs($d.Paragraphs.Count).Range
$lastPar.Find.Execute("<p> <b>NOTE:", $true, $false, $false, $false, $false, $true, 1, $false, "NOTE:", 2)
